# Applies the "Further added models and implemented some textures" edit
# to the PostalFrog "Asset list" workbook.
#
# Summary of the change:
#   - A new "HouseDoor" asset row is appended (row 18).
#   - The "house" asset (row 11) is renamed to "HouseWall", its appearance
#     text is simplified, and its poly count drops from 50 to 25.
#   - Texture resolution (column E) is doubled for most rows (512->1024,
#     1024->2048) as higher-res textures were implemented.
#   - The saved selection moves to H1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New row 18 (HouseDoor asset) and row 11 rename ("house" -> "HouseWall").
#    The string-valued cells are written in the same order the original
#    author typed them so the rebuilt shared-strings table lines up
#    exactly (HouseDoor, Simple wood door, white brick house with window
#    and door, 82 x 204 x 5, HouseWall).
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "HouseDoor"
$ws.Range("C18").Value = "Simple wood door"
$ws.Range("C11").Value = "white brick house with window and door"
$ws.Range("F18").Value = "82 x 204 x 5"
$ws.Range("A11").Value = "HouseWall"

$ws.Range("B18").Value = "Scene Art"
$ws.Range("D18").Value = 30
$ws.Range("E18").Value = 1024
$ws.Range("G18").Value = "low"
$ws.Range("D11").Value = 25

# ---------------------------------------------------------------------
# 2. Texture resolution (column E) doublings across most rows.
# ---------------------------------------------------------------------
$ws.Range("E2").Value = 2048
$ws.Range("E3").Value = 1024
$ws.Range("E4").Value = 1024
$ws.Range("E5").Value = 1024
$ws.Range("E6").Value = 1024
$ws.Range("E7").Value = 1024
$ws.Range("E8").Value = 2048
$ws.Range("E9").Value = 2048
$ws.Range("E10").Value = 2048
$ws.Range("E11").Value = 2048
$ws.Range("E12").Value = 2048
$ws.Range("E13").Value = 1024
$ws.Range("E14").Value = 1024
$ws.Range("E16").Value = 1024
$ws.Range("E17").Value = 2048

# ---------------------------------------------------------------------
# 3. Update the active selection to reflect the post-edit state (H1).
# ---------------------------------------------------------------------
$ws.Range("H1").Select()
